# Updated cryptos list on Tue May 28 23:23:34 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'68.421.71"
$ws.Range("E2").Value = "  -1.46%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.852.27"
$ws.Range("E3").Value = "  -0.88%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'602.10"
$ws.Range("E5").Value = "  -0.29%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'168.59"
$ws.Range("E6").Value = "  -1.19%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "'3.852.36"
$ws.Range("E7").Value = "  -0.90%  "

# Row 8 - USDC
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.18%  "

# Row 9 - XRP
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  -1.12%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -2.23%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +1.09%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -2.60%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +4.02%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'36.98"
$ws.Range("E14").Value = "  -3.41%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'4.497.85"
$ws.Range("E15").Value = "  -0.99%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "'3.855.73"
$ws.Range("E16").Value = "  -0.97%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "'68.536.09"
$ws.Range("E17").Value = "  -1.44%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "'18.51"
$ws.Range("E18").Value = "  -1.29%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -3.33%  "

# Row 20 - was Uniswap, now TRON
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "'0.111"
$ws.Range("E20").Value = "  -1.00%  "

# Row 21 - was TRON, now Uniswap
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'11.19"
$ws.Range("E21").Value = "  +1.14%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'470.66"
$ws.Range("E22").Value = "  -3.84%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "'0.731"
$ws.Range("E23").Value = "  -1.68%  "

# Row 24 - PEPE
$ws.Range("D24").Value = "'0.0000160"
$ws.Range("E24").Value = "  -3.66%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "'83.42"
$ws.Range("E25").Value = "  -2.21%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  -2.94%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "'12.13"
$ws.Range("E27").Value = "  -2.02%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "'10.17"
$ws.Range("E28").Value = "  +0.57%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.12%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "'2.96"
$ws.Range("E30").Value = "  -1.10%  "

# Row 31 - WrappedeETH
$ws.Range("D31").Value = "'4.003.43"
$ws.Range("E31").Value = "  -0.89%  "

# Row 32 - NEARProtocol
$ws.Range("E32").Value = "  -1.86%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "'31.41"
$ws.Range("E33").Value = "  -1.54%  "

# Row 34 - ImmutableX
$ws.Range("E34").Value = "  -4.01%  "

# Row 35 - Aptos
$ws.Range("E35").Value = "  -3.31%  "

# Row 36 - RenzoRestakedETH
$ws.Range("D36").Value = "'3.817.63"
$ws.Range("E36").Value = "  -1.00%  "

# Row 37 - was Hedera, now dogwifhat
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.75"
$ws.Range("E37").Value = "  +10.14%  "

# Row 38 - was dogwifhat, now Hedera
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.104"
$ws.Range("E38").Value = "  -2.53%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -1.90%  "

# Row 40 - Mantle
$ws.Range("E40").Value = "  -1.52%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  -2.85%  "

# Row 42 - FirstDigitalUSD
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.08%  "

# Row 43 - TheGraph
$ws.Range("E43").Value = "  -3.87%  "

# Row 44 - Stacks
$ws.Range("E44").Value = "  -4.63%  "

# Row 45 - Cosmos
$ws.Range("D45").Value = "'8.69"
$ws.Range("E45").Value = "  -0.20%  "

# Row 46 - Bittensor
$ws.Range("D46").Value = "'417.72"
$ws.Range("E46").Value = "  -4.33%  "

# Row 47 - USDe : no change

# Row 48 - was OKB, now FLOKI
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").Value = "'0.000293"
$ws.Range("E48").Value = "  +5.94%  "

# Row 49 - was FLOKI, now OKB
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'46.93"
$ws.Range("E49").Value = "  -2.34%  "

# Row 50 - was EnergySwap, now Monero
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'141.93"
$ws.Range("E50").Value = "  +0.52%  "

# Row 51 - was Monero, now EnergySwap
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'26.07"
$ws.Range("E51").Value = "  +3.62%  "
